# Auto-generated Excel COM-interop script
# Applies the numeric updates described in the commit diff to the
# Golem_Profits workbook (per-sheet leve profit recalculations).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 599.7646999999999
$ws.Range("I28").Value = 616.4
$ws.Range("K28").Value = 616.4
$ws.Range("M28").Value = -131.4

$ws.Range("H53").Value = 212.8
$ws.Range("I53").Value = 220.66667
$ws.Range("J53").Value = 201
$ws.Range("K53").Value = 220.66667
$ws.Range("L53").Value = 201
$ws.Range("M53").Value = 416.33333
$ws.Range("N53").Value = -1475

$ws.Range("H80").Value = 875.0909
$ws.Range("I80").Value = 943.6
$ws.Range("J80").Value = 818
$ws.Range("K80").Value = 2830.8
$ws.Range("L80").Value = 2454
$ws.Range("M80").Value = -1832.8
$ws.Range("N80").Value = -4450

$ws.Range("H83").Value = 875.0909
$ws.Range("I83").Value = 943.6
$ws.Range("J83").Value = 818
$ws.Range("K83").Value = 8492.4
$ws.Range("L83").Value = 7362
$ws.Range("M83").Value = -3500.4
$ws.Range("N83").Value = -17346

$ws.Range("H101").Value = 403.5
$ws.Range("I101").Value = 403.5
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1210.5
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 411.5
$ws.Range("N101").ClearContents()

$ws.Range("H105").Value = 49333.332
$ws.Range("J105").Value = 49333.332
$ws.Range("L105").Value = 49333.332
$ws.Range("N105").Value = -56321.332

$ws.Range("H132").Value = 1005
$ws.Range("I132").Value = 1005
$ws.Range("K132").Value = 3015
$ws.Range("M132").Value = -485


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 288
$ws.Range("I2").Value = 297.05554
$ws.Range("K2").Value = 297.05554
$ws.Range("M2").Value = -184.05554

$ws.Range("H45").Value = 1424.4897
$ws.Range("I45").Value = 1412.5
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1412.5
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1035.5
$ws.Range("N45").Value = -2754

$ws.Range("H61").Value = 499.5
$ws.Range("I61").Value = 499.5
$ws.Range("K61").Value = 499.5
$ws.Range("M61").Value = -287.5

$ws.Range("H74").Value = 2822.4
$ws.Range("I74").Value = 2303
$ws.Range("K74").Value = 2303
$ws.Range("M74").Value = -1429

$ws.Range("H77").Value = 2822.4
$ws.Range("I77").Value = 2303
$ws.Range("K77").Value = 11515
$ws.Range("M77").Value = -7147

$ws.Range("H110").Value = 937.6
$ws.Range("I110").Value = 1099.6666
$ws.Range("J110").Value = 694.5
$ws.Range("K110").Value = 1099.6666
$ws.Range("L110").Value = 694.5
$ws.Range("M110").Value = 945.3334
$ws.Range("N110").Value = -4784.5

$ws.Range("H116").Value = 288
$ws.Range("I116").Value = 297.05554
$ws.Range("K116").Value = 297.05554
$ws.Range("M116").Value = 1996.94446

$ws.Range("H136").Value = 499.5
$ws.Range("I136").Value = 499.5
$ws.Range("K136").Value = 1498.5
$ws.Range("M136").Value = 1051.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 288
$ws.Range("I3").Value = 297.05554
$ws.Range("K3").Value = 297.05554
$ws.Range("M3").Value = -183.05554

$ws.Range("H134").Value = 368.5
$ws.Range("I134").Value = 368.5
$ws.Range("K134").Value = 1105.5
$ws.Range("M134").Value = 1429.5


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 851.6667
$ws.Range("J16").Value = 777.5
$ws.Range("L16").Value = 777.5
$ws.Range("N16").Value = -1351.5

$ws.Range("H22").Value = 749.93335
$ws.Range("I22").Value = 769
$ws.Range("J22").Value = 626
$ws.Range("K22").Value = 769
$ws.Range("L22").Value = 626
$ws.Range("M22").Value = -419
$ws.Range("N22").Value = -1326

$ws.Range("H113").Value = 851.6667
$ws.Range("J113").Value = 777.5
$ws.Range("L113").Value = 777.5
$ws.Range("N113").Value = -5117.5

$ws.Range("H122").Value = 1081.7142
$ws.Range("I122").Value = 1034.4
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3103.2
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -653.2000000000003
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 2238.25
$ws.Range("I132").Value = 1984.3334
$ws.Range("K132").Value = 5953.0002
$ws.Range("M132").Value = -3423.0002

$ws.Range("H134").Value = 1333.3334
$ws.Range("I134").Value = 1333.3334
$ws.Range("K134").Value = 4000.0002
$ws.Range("M134").Value = -1465.0002


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 11875.583
$ws.Range("I25").Value = 2501.1667
$ws.Range("J25").Value = 21250
$ws.Range("K25").Value = 7503.500100000001
$ws.Range("L25").Value = 63750
$ws.Range("M25").Value = -7334.500100000001
$ws.Range("N25").Value = -64088

$ws.Range("H30").Value = 11875.583
$ws.Range("I30").Value = 2501.1667
$ws.Range("J30").Value = 21250
$ws.Range("K30").Value = 7503.500100000001
$ws.Range("L30").Value = 63750
$ws.Range("M30").Value = -7401.500100000001
$ws.Range("N30").Value = -63954

$ws.Range("H55").Value = 2116.6667
$ws.Range("I55").Value = 1466.6666
$ws.Range("J55").Value = 4066.6667
$ws.Range("K55").Value = 4399.9998
$ws.Range("L55").Value = 12200.0001
$ws.Range("M55").Value = -4222.9998
$ws.Range("N55").Value = -12554.0001

$ws.Range("H140").Value = 430
$ws.Range("I140").Value = 430
$ws.Range("K140").Value = 1290
$ws.Range("M140").Value = 3890


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 502500
$ws.Range("J21").Value = 502500
$ws.Range("L21").Value = 502500
$ws.Range("N21").Value = -502846

$ws.Range("H30").Value = 502500
$ws.Range("J30").Value = 502500
$ws.Range("L30").Value = 502500
$ws.Range("N30").Value = -502710

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H97").Value = 487.375
$ws.Range("I97").Value = 467.7143
$ws.Range("J97").Value = 625
$ws.Range("K97").Value = 467.7143
$ws.Range("L97").Value = 625
$ws.Range("M97").Value = 28.28570000000002
$ws.Range("N97").Value = -1617

$ws.Range("H132").Value = 499
$ws.Range("I132").Value = 499
$ws.Range("K132").Value = 1497
$ws.Range("M132").Value = 1033


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2801.3635
$ws.Range("I22").Value = 1571.5714
$ws.Range("J22").Value = 4953.5
$ws.Range("K22").Value = 1571.5714
$ws.Range("L22").Value = 4953.5
$ws.Range("M22").Value = -1276.5714
$ws.Range("N22").Value = -5543.5

$ws.Range("H27").Value = 2801.3635
$ws.Range("I27").Value = 1571.5714
$ws.Range("J27").Value = 4953.5
$ws.Range("K27").Value = 1571.5714
$ws.Range("L27").Value = 4953.5
$ws.Range("M27").Value = -1464.5714
$ws.Range("N27").Value = -5167.5

$ws.Range("H40").Value = 851834
$ws.Range("I40").Value = 22200.8
$ws.Range("J40").Value = 5000000
$ws.Range("K40").Value = 22200.8
$ws.Range("L40").Value = 5000000
$ws.Range("M40").Value = -22064.8
$ws.Range("N40").Value = -5000272

$ws.Range("H46").Value = 407996.6
$ws.Range("I46").Value = 2000000
$ws.Range("K46").Value = 2000000
$ws.Range("M46").Value = -1999812


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H51").Value = 20077
$ws.Range("J51").Value = 20077
$ws.Range("L51").Value = 20077
$ws.Range("N51").Value = -21097

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H124").Value = 36249.75
$ws.Range("J124").Value = 43333
$ws.Range("L124").Value = 43333
$ws.Range("N124").Value = -53153

$ws.Range("H136").Value = 1022.5
$ws.Range("I136").Value = 842.8182
$ws.Range("K136").Value = 2528.4546
$ws.Range("M136").Value = 21.54539999999997

